$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 and IF headers, matching the style of
# the existing header cells (bold, thin border, centered/top aligned). ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-15: new values for columns I (I0) and J (IF) ---
$data = @{
    2  = @(5, 5)
    3  = @(6, 6)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 5)
    7  = @(1, 5)
    8  = @(1, 5)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 4)
    13 = @(1, 3)
    14 = @(3, 4)
    15 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

$excel.DisplayAlerts = $false
Write-Host "I0 and IF columns added"
